$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table shifted up by one row: the original row 2 (date 39400 / 2007)
# was dropped, and every subsequent row moved up to take its place. This also
# drops the final (former row 19) data point, shrinking the table from
# A1:E19 down to A1:E18. Deleting row 2 performs exactly this shift.
$ws.Rows(2).Delete()

# The forecast column (E) values were recalculated for every remaining data
# row (rows 2-18 after the shift above).
$ws.Range("E2").Value = 0.4944284391569687
$ws.Range("E3").Value = -0.4782015746048418
$ws.Range("E4").Value = 1.324233212457782
$ws.Range("E5").Value = 0.7478380109886329
$ws.Range("E6").Value = -0.2445716668737163
$ws.Range("E7").Value = -0.2617076051026235
$ws.Range("E8").Value = -0.100009932057743
$ws.Range("E9").Value = 0.3000376062062493
$ws.Range("E10").Value = 0.1740313431290996
$ws.Range("E11").Value = 0.3390041783450259
$ws.Range("E12").Value = 0.2210188332817387
$ws.Range("E13").Value = -0.09571633453315798
$ws.Range("E14").Value = -1.49562970548649
$ws.Range("E15").Value = -0.1048501255800471
$ws.Range("E16").Value = 0.9692952624595019
$ws.Range("E17").Value = 0.1544084105021826
$ws.Range("E18").Value = 0.3997355152047577
